$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 14 (BMP280 component) - everything below shifts up
$ws.Rows.Item(14).Delete()

# The SCD41-D-R2 row (now row 15 after the shift) had its CO2 calibration
# range comment updated from "400-2000 ppm" to "400-5000 ppm"
$ws.Range("N15").Value = "400-5000 ppm, Stäng av autocalibrate"

# Update selection to reflect where the editor last clicked
$ws.Range("C29").Select()
